# Auto-generated edit script: apply numeric cell updates per sheet/row
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3615.24
$ws.Range("I32").Value = 3417.8145
$ws.Range("J32").Value = 9998.666999999999
$ws.Range("K32").Value = 3417.8145
$ws.Range("L32").Value = 9998.666999999999
$ws.Range("M32").Value = -3130.8145
$ws.Range("N32").Value = -10572.667

$ws.Range("H45").Value = 2782.3333
$ws.Range("I45").Value = 3234.4285
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 3234.4285
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -2857.4285
$ws.Range("N45").Value = -1954

$ws.Range("H102").Value = 12822107
$ws.Range("I102").Value = 16667889
$ws.Range("J102").Value = 2833.3333
$ws.Range("K102").Value = 16667889
$ws.Range("L102").Value = 2833.3333
$ws.Range("M102").Value = -16666267

$ws.Range("H122").Value = 3780
$ws.Range("I122").Value = 3385.3333
$ws.Range("J122").Value = 4964
$ws.Range("K122").Value = 10155.9999
$ws.Range("L122").Value = 14892
$ws.Range("M122").Value = -7705.999899999999
$ws.Range("N122").Value = -19792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 861.1667
$ws.Range("I80").Value = 348.75
$ws.Range("J80").Value = 1117.375
$ws.Range("K80").Value = 348.75
$ws.Range("L80").Value = 1117.375
$ws.Range("M80").Value = 649.25

$ws.Range("H83").Value = 861.1667
$ws.Range("I83").Value = 348.75
$ws.Range("J83").Value = 1117.375
$ws.Range("K83").Value = 1743.75
$ws.Range("L83").Value = 5586.875
$ws.Range("M83").Value = 3248.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H108").Value = 30927
$ws.Range("I108").Value = 20621
$ws.Range("J108").Value = 33503.5
$ws.Range("K108").Value = 20621
$ws.Range("L108").Value = 33503.5
$ws.Range("M108").Value = -16781
$ws.Range("N108").Value = -41183.5

$ws.Range("H131").Value = 11930.571
$ws.Range("I131").Value = 4296
$ws.Range("J131").Value = 17656.5
$ws.Range("K131").Value = 4296
$ws.Range("L131").Value = 17656.5
$ws.Range("M131").Value = 744
$ws.Range("N131").Value = -27736.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 49.5
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 71
$ws.Range("K2").Value = 168
$ws.Range("L2").Value = 426
$ws.Range("M2").Value = -55
$ws.Range("N2").Value = -652

$ws.Range("H4").Value = 5180462
$ws.Range("I4").Value = 4000127.8
$ws.Range("J4").Value = 6360796
$ws.Range("K4").Value = 12000383.4
$ws.Range("L4").Value = 19082388
$ws.Range("M4").Value = -12000271.4
$ws.Range("N4").Value = -19082612

$ws.Range("H40").Value = 205.86667
$ws.Range("I40").Value = 76.44444
$ws.Range("J40").Value = 400
$ws.Range("K40").Value = 305.77776
$ws.Range("L40").Value = 1600
$ws.Range("M40").Value = -236.77776

$ws.Range("H107").Value = 4569.2085
$ws.Range("I107").Value = 429.83334
$ws.Range("J107").Value = 5949
$ws.Range("K107").Value = 1289.50002
$ws.Range("L107").Value = 17847
$ws.Range("M107").Value = 630.4999800000001
$ws.Range("N107").Value = -21687

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19995
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 19995
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 19995
$ws.Range("N15").Value = -20571

$ws.Range("H52").Value = 20599.8
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 20599.8
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 20599.8
$ws.Range("N52").Value = -21117.8

$ws.Range("H59").Value = 10000
$ws.Range("I59").Value = 10000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -9417
$ws.Range("N59").ClearContents()

$ws.Range("H81").Value = 19995
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 19995
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 19995
$ws.Range("N81").Value = -21991

$ws.Range("H84").Value = 19995
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 19995
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 59985
$ws.Range("N84").Value = -69969

$ws.Range("H132").Value = 3572.3044
$ws.Range("I132").Value = 3929.25
$ws.Range("J132").Value = 3182.9092
$ws.Range("K132").Value = 11787.75
$ws.Range("L132").Value = 9548.7276
$ws.Range("M132").Value = -9257.75
$ws.Range("N132").Value = -14608.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H46").Value = 3923.4707
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 4118.6875
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 4118.6875
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -4494.6875

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()

$ws.Range("H93").Value = 997.58826
$ws.Range("I93").Value = 959.9375
$ws.Range("J93").Value = 1600
$ws.Range("K93").Value = 959.9375
$ws.Range("L93").Value = 1600
$ws.Range("M93").Value = 288.0625

$ws.Range("H122").Value = 20835114
$ws.Range("I122").Value = 31251544
$ws.Range("J122").Value = 2251
$ws.Range("K122").Value = 93754632
$ws.Range("L122").Value = 6753
$ws.Range("M122").Value = -93752182
$ws.Range("N122").Value = -11653

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4861.4614
$ws.Range("I96").Value = 4319.9
$ws.Range("J96").Value = 6666.6665
$ws.Range("K96").Value = 4319.9
$ws.Range("L96").Value = 6666.6665
$ws.Range("M96").Value = -2946.9
$ws.Range("N96").Value = -9412.666499999999

$ws.Range("H100").Value = 399.375
$ws.Range("I100").Value = 387.8
$ws.Range("J100").Value = 418.66666
$ws.Range("K100").Value = 775.6
$ws.Range("L100").Value = 837.33332
$ws.Range("M100").Value = -234.6
$ws.Range("N100").Value = -1919.33332

$ws.Range("H132").Value = 1712.8
$ws.Range("I132").Value = 1605.8611
$ws.Range("J132").Value = 2140.5557
$ws.Range("K132").Value = 4817.5833
$ws.Range("L132").Value = 6421.6671
$ws.Range("M132").Value = -2287.5833
$ws.Range("N132").Value = -11481.6671

$ws.Range("H136").Value = 1513.6
$ws.Range("I136").Value = 1323.8889
$ws.Range("J136").Value = 1668.8182
$ws.Range("K136").Value = 3971.6667
$ws.Range("L136").Value = 5006.4546
$ws.Range("M136").Value = -1421.6667
